$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.366.05"
$ws.Range("E2").Value = "  +0.15%  "

$ws.Range("D3").Value = "2.067.00"
$ws.Range("E3").Value = "  +0.35%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'235.33"
$ws.Range("E5").Value = "  -0.10%  "

$ws.Range("E6").Value = "  +2.16%  "

$ws.Range("D8").Value = "'57.59"
$ws.Range("E8").Value = "  -0.80%  "

$ws.Range("D9").Value = "'0.395"
$ws.Range("E9").Value = "  +3.39%  "

$ws.Range("D10").Value = "'0.0773"
$ws.Range("E10").Value = "  +1.38%  "

$ws.Range("D11").Value = "'0.103"
$ws.Range("E11").Value = "  +0.79%  "

$ws.Range("D12").Value = "2.370.01"
$ws.Range("E12").Value = "  +0.27%  "

$ws.Range("D13").Value = "'14.39"
$ws.Range("E13").Value = "  -0.80%  "

$ws.Range("D14").Value = "'20.72"
$ws.Range("E14").Value = "  -0.92%  "

$ws.Range("E15").Value = "  -0.12%  "

$ws.Range("D16").Value = "'5.19"
$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("D17").Value = "2.065.43"
$ws.Range("E17").Value = "  +0.21%  "

$ws.Range("D18").Value = "37.302.57"
$ws.Range("E18").Value = "  -0.69%  "

$ws.Range("D19").Value = "'6.19"
$ws.Range("E19").Value = "  -0.47%  "

$ws.Range("D20").Value = "'69.60"
$ws.Range("E20").Value = "  +0.82%  "

$ws.Range("D21").Value = "0.0₃0817"
$ws.Range("E21").Value = "  +0.25%  "

$ws.Range("D22").Value = "'226.92"
$ws.Range("E22").Value = "  +0.40%  "

$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("E24").Value = "  +2.21%  "

$ws.Range("E25").Value = "  -1.55%  "

$ws.Range("D26").Value = "'166.84"
$ws.Range("E26").Value = "  +1.48%  "

$ws.Range("D27").Value = "'8.91"
$ws.Range("E27").Value = "  +0.41%  "

$ws.Range("D28").Value = "'1.41"
$ws.Range("E28").Value = "  -5.19%  "

$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "'0.127"
$ws.Range("E29").Value = "  -0.76%  "

$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'19.12"
$ws.Range("E30").Value = "  -0.47%  "

$ws.Range("D31").Value = "'0.118"
$ws.Range("E31").Value = "  -0.85%  "

$ws.Range("E32").Value = "  +1.05%  "

$ws.Range("E33").Value = "  -0.93%  "

$ws.Range("D34").Value = "'4.55"
$ws.Range("E34").Value = "  +1.52%  "

$ws.Range("D35").Value = "'2.48"
$ws.Range("E35").Value = "  -3.20%  "

$ws.Range("E36").Value = "  +0.15%  "

$ws.Range("D37").Value = "'3.34"
$ws.Range("E37").Value = "  -2.74%  "

$ws.Range("E38").Value = "  +0.16%  "

$ws.Range("D39").Value = "'5.63"
$ws.Range("E39").Value = "  -4.56%  "

$ws.Range("E40").Value = "  -0.80%  "

$ws.Range("D41").Value = "'0.0959"
$ws.Range("E41").Value = "  -2.66%  "

$ws.Range("D42").Value = "'97.65"
$ws.Range("E42").Value = "  +0.79%  "

$ws.Range("D43").Value = "1.480.19"
$ws.Range("E43").Value = "  +0.39%  "

$ws.Range("D44").Value = "'0.0212"
$ws.Range("E44").Value = "  +0.87%  "

$ws.Range("D45").Value = "'1.16"
$ws.Range("E45").Value = "  +0.15%  "

$ws.Range("D46").Value = "'4.06"
$ws.Range("E46").Value = "  -10.65%  "

$ws.Range("D47").Value = "'1.02"
$ws.Range("E47").Value = "  +0.10%  "

$ws.Range("D48").Value = "'15.29"
$ws.Range("E48").Value = "  -3.92%  "

$ws.Range("D49").Value = "'7.22"
$ws.Range("E49").Value = "  +0.16%  "

$ws.Range("D50").Value = "'2.96"
$ws.Range("E50").Value = "  +0.77%  "

$ws.Range("D51").Value = "2.257.01"
$ws.Range("E51").Value = "  +0.27%  "
